$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 23349
$ws.Range("D2").Value = 34099795
$ws.Range("C3").Value = 58831
$ws.Range("D3").Value = 87111745
$ws.Range("C4").Value = 19845
$ws.Range("D4").Value = 29627823
$ws.Range("C5").Value = 5246
$ws.Range("D5").Value = 7853739
$ws.Range("C6").Value = 1028
$ws.Range("D6").Value = 1540391
$ws.Range("C7").Value = 65
$ws.Range("D7").Value = 97500
$ws.Range("C8").Value = 5
$ws.Range("D8").Value = 7500
$ws.Range("C10").Value = 24995
$ws.Range("D10").Value = 34563220
$ws.Range("C11").Value = 6124
$ws.Range("D11").Value = 8962808
$ws.Range("C12").Value = 17266
$ws.Range("D12").Value = 25550309
$ws.Range("C13").Value = 5352
$ws.Range("D13").Value = 7996759
$ws.Range("C14").Value = 1260
$ws.Range("D14").Value = 1885995
$ws.Range("C15").Value = 229
$ws.Range("D15").Value = 340766
$ws.Range("C16").Value = 18
$ws.Range("D16").Value = 27000
$ws.Range("C17").Value = 6290
$ws.Range("D17").Value = 8525015
$ws.Range("C18").Value = 8565
$ws.Range("D18").Value = 12491005
$ws.Range("C19").Value = 21171
$ws.Range("D19").Value = 31354655
$ws.Range("C20").Value = 6691
$ws.Range("D20").Value = 10003916
$ws.Range("C21").Value = 1600
$ws.Range("D21").Value = 2395670
$ws.Range("C22").Value = 237
$ws.Range("D22").Value = 355124
$ws.Range("C24").Value = 7332
$ws.Range("D24").Value = 10024668
$ws.Range("C25").Value = 4843
$ws.Range("D25").Value = 7069794
$ws.Range("C26").Value = 14902
$ws.Range("D26").Value = 22051807
$ws.Range("C27").Value = 5042
$ws.Range("D27").Value = 7540372
$ws.Range("C28").Value = 1207
$ws.Range("D28").Value = 1809991
$ws.Range("C29").Value = 177
$ws.Range("D29").Value = 265500
$ws.Range("C31").Value = 5176
$ws.Range("D31").Value = 6956870
$ws.Range("C32").Value = 1707
$ws.Range("D32").Value = 2463796
$ws.Range("C33").Value = 4519
$ws.Range("D33").Value = 6642540
$ws.Range("C34").Value = 1827
$ws.Range("D34").Value = 2720154
$ws.Range("C35").Value = 468
$ws.Range("D35").Value = 698541
$ws.Range("C36").Value = 88
$ws.Range("D36").Value = 132000
$ws.Range("C38").Value = 1154
$ws.Range("D38").Value = 1569940
$ws.Range("C39").Value = 10864
$ws.Range("D39").Value = 15855628
$ws.Range("C40").Value = 33528
$ws.Range("D40").Value = 49597265
$ws.Range("C41").Value = 12360
$ws.Range("D41").Value = 18475351
$ws.Range("C42").Value = 3414
$ws.Range("D42").Value = 5112941
$ws.Range("C43").Value = 591
$ws.Range("D43").Value = 885436
$ws.Range("C46").Value = 10259
$ws.Range("D46").Value = 13983390
$ws.Range("C47").Value = 969
$ws.Range("D47").Value = 1403195
$ws.Range("C48").Value = 3641
$ws.Range("D48").Value = 5369139
$ws.Range("C49").Value = 1374
$ws.Range("D49").Value = 2053964
$ws.Range("C50").Value = 421
$ws.Range("D50").Value = 629000
$ws.Range("C51").Value = 87
$ws.Range("D51").Value = 130500
$ws.Range("C52").Value = 2319
$ws.Range("D52").Value = 3228118
$ws.Range("C53").Value = 346
$ws.Range("D53").Value = 502784
$ws.Range("C54").Value = 927
$ws.Range("D54").Value = 1373977
$ws.Range("C55").Value = 375
$ws.Range("D55").Value = 560476
$ws.Range("C56").Value = 127
$ws.Range("D56").Value = 190378
$ws.Range("C58").Value = 436
$ws.Range("D58").Value = 622722
$ws.Range("C59").Value = 9883
$ws.Range("D59").Value = 14363165
$ws.Range("C60").Value = 30081
$ws.Range("D60").Value = 44401794
$ws.Range("C61").Value = 10393
$ws.Range("D61").Value = 15539264
$ws.Range("C62").Value = 2881
$ws.Range("D62").Value = 4312068
$ws.Range("C63").Value = 502
$ws.Range("D63").Value = 752639
$ws.Range("C64").Value = 40
$ws.Range("D64").Value = 60000
$ws.Range("C66").Value = 9782
$ws.Range("D66").Value = 13114433
$ws.Range("C67").Value = 2677
$ws.Range("D67").Value = 3909805
$ws.Range("C68").Value = 7279
$ws.Range("D68").Value = 10742239
$ws.Range("C69").Value = 2576
$ws.Range("D69").Value = 3849272
$ws.Range("C70").Value = 840
$ws.Range("D70").Value = 1258049
$ws.Range("C71").Value = 169
$ws.Range("D71").Value = 252112
$ws.Range("C73").Value = 2788
$ws.Range("D73").Value = 3796155
$ws.Range("C74").Value = 849
$ws.Range("D74").Value = 1231250
$ws.Range("C75").Value = 2950
$ws.Range("D75").Value = 4360946
$ws.Range("C76").Value = 1158
$ws.Range("D76").Value = 1733548
$ws.Range("C77").Value = 404
$ws.Range("D77").Value = 606000
$ws.Range("C78").Value = 83
$ws.Range("D78").Value = 124069
$ws.Range("C80").Value = 1736
$ws.Range("D80").Value = 2335993
$ws.Range("C81").Value = 300
$ws.Range("D81").Value = 446189
$ws.Range("C82").Value = 101
$ws.Range("D82").Value = 151110
$ws.Range("C86").Value = 6946
$ws.Range("D86").Value = 10162402
$ws.Range("C87").Value = 19917
$ws.Range("D87").Value = 29477634
$ws.Range("C88").Value = 6536
$ws.Range("D88").Value = 9770766
$ws.Range("C89").Value = 1729
$ws.Range("D89").Value = 2589155
$ws.Range("C90").Value = 275
$ws.Range("D90").Value = 412310
$ws.Range("C91").Value = 22
$ws.Range("D91").Value = 33000
$ws.Range("C93").Value = 6226
$ws.Range("D93").Value = 8399938
$ws.Range("C94").Value = 19080
$ws.Range("D94").Value = 27720752
$ws.Range("C95").Value = 44190
$ws.Range("D95").Value = 65217034
$ws.Range("C96").Value = 14121
$ws.Range("D96").Value = 21086660
$ws.Range("C97").Value = 3751
$ws.Range("D97").Value = 5615102
$ws.Range("C98").Value = 640
$ws.Range("D98").Value = 958362
$ws.Range("C101").Value = 16286
$ws.Range("D101").Value = 22147311
$ws.Range("C102").Value = 21810
$ws.Range("D102").Value = 31726482
$ws.Range("C103").Value = 49247
$ws.Range("D103").Value = 72590374
$ws.Range("C104").Value = 15336
$ws.Range("D104").Value = 22880870
$ws.Range("C105").Value = 3926
$ws.Range("D105").Value = 5865922
$ws.Range("C106").Value = 633
$ws.Range("D106").Value = 946554
$ws.Range("C109").Value = 19315
$ws.Range("D109").Value = 26084628
$ws.Range("C110").Value = 8509
$ws.Range("D110").Value = 12435326
$ws.Range("C111").Value = 21981
$ws.Range("D111").Value = 32553924
$ws.Range("C112").Value = 7580
$ws.Range("D112").Value = 11318049
$ws.Range("C113").Value = 1841
$ws.Range("D113").Value = 2755513
$ws.Range("C114").Value = 263
$ws.Range("D114").Value = 392728
$ws.Range("C117").Value = 6912
$ws.Range("D117").Value = 9421084
$ws.Range("C118").Value = 21214
$ws.Range("D118").Value = 30857240
$ws.Range("C119").Value = 52208
$ws.Range("D119").Value = 77031661
$ws.Range("C120").Value = 15697
$ws.Range("D120").Value = 23446980
$ws.Range("C121").Value = 3909
$ws.Range("D121").Value = 5849208
$ws.Range("C122").Value = 792
$ws.Range("D122").Value = 1186212
$ws.Range("C125").Value = 18043
$ws.Range("D125").Value = 24819156
$ws.Range("C126").Value = 28978
$ws.Range("D126").Value = 42452372
$ws.Range("C127").Value = 87246
$ws.Range("D127").Value = 129282088
$ws.Range("C128").Value = 38642
$ws.Range("D128").Value = 57767798
$ws.Range("C129").Value = 12031
$ws.Range("D129").Value = 18021724
$ws.Range("C130").Value = 2420
$ws.Range("D130").Value = 3625602
$ws.Range("C131").Value = 130
$ws.Range("D131").Value = 194212
$ws.Range("C134").Value = 28569
$ws.Range("D134").Value = 39879744
